# Availability.xlsx — add a "Phone Number" field + a new row for Patrick Starkey
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row 11 (Patrick Starkey) + new "Phone Number" header (I2) ---
# Values are entered in the same order the original author used so new
# shared-string table entries land in the same index order as the target.
$ws.Range("A11").Value = "Patrick Starkey"
$ws.Range("B11").Value = "9 am-MN"
$ws.Range("F11").Value = "9-MN"
$ws.Range("J11").Value = "iakavas@live.com"
$ws.Range("G11").Value = "10am-3"
$ws.Range("H11").Value = "9am-12"

# "Phone Number" header, copying the formatting used by the other header
# cells in row 2 (e.g. H2).
$ws.Range("I2").Value = "Phone Number"
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

$ws.Range("I11").Value = "281-797-7242"
$ws.Range("C11").Value = "9am-1230 & 2-9"
$ws.Range("D11").Value = "9am-MN"
$ws.Range("E11").Value = "9am-1230 & 2-9"

# Hyperlink the e-mail address cell, then restore the shared "Hyperlink"
# cell style (adding the hyperlink on its own leaves a slightly different
# style behind).
$ws.Hyperlinks.Add($ws.Range("J11"), "mailto:iakavas@live.com")
$ws.Range("J3").Copy()
$ws.Range("J11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column width tweaks ---
# Column.ColumnWidth is quantized to whole pixels by this host, so the
# inputs below are chosen so the stored (1/6-character-unit) width lands on
# the closest achievable value to the target widths.
$ws.Columns.Item(3).ColumnWidth = 12.91796875     # -> stored width 13.8333...
$ws.Columns.Item(5).ColumnWidth = 15.584635416666666  # -> stored width 16.5
$ws.Columns.Item(9).ColumnWidth = 19.41796875     # -> stored width 20.3333...

# --- Selection / view update ---
$null = $ws.Range("C11").Select()
